$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035073802085231
$ws.Cells.Item(2, 4).Value = 1.044883814922542
$ws.Cells.Item(2, 5).Value = 1.053160801652799
$ws.Cells.Item(2, 6).Value = 1.058778293861671
$ws.Cells.Item(2, 9).Value = 1.042831549444113
$ws.Cells.Item(2, 10).Value = 1.040189568209835
$ws.Cells.Item(2, 11).Value = 1.047653611860067
$ws.Cells.Item(2, 12).Value = 1.055907536164126
$ws.Cells.Item(2, 13).Value = 1.061509598669011
$ws.Cells.Item(2, 14).Value = 1.017369357015482

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.03582037655637
$ws.Cells.Item(3, 4).Value = 1.045472992837573
$ws.Cells.Item(3, 5).Value = 1.053970331342616
$ws.Cells.Item(3, 6).Value = 1.059538130331408
$ws.Cells.Item(3, 9).Value = 1.043016063171804
$ws.Cells.Item(3, 10).Value = 1.040580564541448
$ws.Cells.Item(3, 11).Value = 1.048054697361961
$ws.Cells.Item(3, 12).Value = 1.056530064234971
$ws.Cells.Item(3, 13).Value = 1.062083672563085
$ws.Cells.Item(3, 14).Value = 1.017499070364507

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036304033056783
$ws.Cells.Item(4, 4).Value = 1.045854727069718
$ws.Cells.Item(4, 5).Value = 1.054495528983501
$ws.Cells.Item(4, 6).Value = 1.060030805465144
$ws.Cells.Item(4, 9).Value = 1.043134439233762
$ws.Cells.Item(4, 10).Value = 1.040833434382065
$ws.Cells.Item(4, 11).Value = 1.04831402299505
$ws.Cells.Item(4, 12).Value = 1.056933602957422
$ws.Cells.Item(4, 13).Value = 1.06245548010675
$ws.Cells.Item(4, 14).Value = 1.017582947521202

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.036507497472811
$ws.Cells.Item(5, 4).Value = 1.04601532511971
$ws.Cells.Item(5, 5).Value = 1.054716649755097
$ws.Cells.Item(5, 6).Value = 1.060238165750552
$ws.Cells.Item(5, 9).Value = 1.043183960247539
$ws.Cells.Item(5, 10).Value = 1.040939708215503
$ws.Cells.Item(5, 11).Value = 1.04842299336856
$ws.Cells.Item(5, 12).Value = 1.057103421562356
$ws.Cells.Item(5, 13).Value = 1.062611868429804
$ws.Cells.Item(5, 14).Value = 1.017618195601563

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036541667917029
$ws.Cells.Item(6, 4).Value = 1.046042297042246
$ws.Cells.Item(6, 5).Value = 1.05475379605817
$ws.Cells.Item(6, 6).Value = 1.060272996445564
$ws.Cells.Item(6, 9).Value = 1.043192260692864
$ws.Cells.Item(6, 10).Value = 1.040957550107411
$ws.Cells.Item(6, 11).Value = 1.048441286993739
$ws.Cells.Item(6, 12).Value = 1.057131944824579
$ws.Cells.Item(6, 13).Value = 1.062638131382562
$ws.Cells.Item(6, 14).Value = 1.017624113082344

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.036306751229352
$ws.Cells.Item(7, 4).Value = 1.045856872531516
$ws.Cells.Item(7, 5).Value = 1.054498482325224
$ws.Cells.Item(7, 6).Value = 1.060033575284457
$ws.Cells.Item(7, 9).Value = 1.043135101896951
$ws.Cells.Item(7, 10).Value = 1.040834854548147
$ws.Cells.Item(7, 11).Value = 1.048315479261224
$ws.Cells.Item(7, 12).Value = 1.056935871412556
$ws.Cells.Item(7, 13).Value = 1.062457569461355
$ws.Cells.Item(7, 14).Value = 1.017583418562876

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.035325990641395
$ws.Cells.Item(8, 4).Value = 1.045082826343351
$ws.Cells.Item(8, 5).Value = 1.053434099417938
$ws.Cells.Item(8, 6).Value = 1.05903487369748
$ws.Cells.Item(8, 9).Value = 1.042894116691068
$ws.Cells.Item(8, 10).Value = 1.040321733578143
$ws.Cells.Item(8, 11).Value = 1.047789201877612
$ws.Cells.Item(8, 12).Value = 1.056117771988899
$ws.Cells.Item(8, 13).Value = 1.061703537475448
$ws.Cells.Item(8, 14).Value = 1.017413205580868

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.033602238235742
$ws.Cells.Item(9, 4).Value = 1.043722748000304
$ws.Cells.Item(9, 5).Value = 1.05156916881499
$ws.Cells.Item(9, 6).Value = 1.057282866151029
$ws.Cells.Item(9, 9).Value = 1.042461723100217
$ws.Cells.Item(9, 10).Value = 1.039416603675656
$ws.Cells.Item(9, 11).Value = 1.046860337290418
$ws.Cells.Item(9, 12).Value = 1.054681780631479
$ws.Cells.Item(9, 13).Value = 1.060377544522397
$ws.Cells.Item(9, 14).Value = 1.017112859446383

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032456192697359
$ws.Cells.Item(10, 4).Value = 1.042818760289379
$ws.Cells.Item(10, 5).Value = 1.050333170513229
$ws.Cells.Item(10, 6).Value = 1.056120256777151
$ws.Cells.Item(10, 9).Value = 1.042168305549472
$ws.Cells.Item(10, 10).Value = 1.038812626060891
$ws.Cells.Item(10, 11).Value = 1.046240171842428
$ws.Cells.Item(10, 12).Value = 1.053728328041574
$ws.Cells.Item(10, 13).Value = 1.059495475410758
$ws.Cells.Item(10, 14).Value = 1.016912380913714

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.031960707386233
$ws.Cells.Item(11, 4).Value = 1.042427995986613
$ws.Cells.Item(11, 5).Value = 1.049799725907993
$ws.Cells.Item(11, 6).Value = 1.055618141370311
$ws.Cells.Item(11, 9).Value = 1.042040041062766
$ws.Cells.Item(11, 10).Value = 1.038550981281008
$ws.Cells.Item(11, 11).Value = 1.045971432330042
$ws.Cells.Item(11, 12).Value = 1.05331641431974
$ws.Cells.Item(11, 13).Value = 1.059114009009518
$ws.Cells.Item(11, 14).Value = 1.016825518502492

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.031776778232341
$ws.Cells.Item(12, 4).Value = 1.042282951136235
$ws.Cells.Item(12, 5).Value = 1.049601846118984
$ws.Cells.Item(12, 6).Value = 1.055431831224807
$ws.Cells.Item(12, 9).Value = 1.041992216618122
$ws.Cells.Item(12, 10).Value = 1.038453778289765
$ws.Cells.Item(12, 11).Value = 1.045871581328907
$ws.Cells.Item(12, 12).Value = 1.053163553926309
$ws.Cells.Item(12, 13).Value = 1.058972388637944
$ws.Cells.Item(12, 14).Value = 1.016793246307063

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.031816226372366
$ws.Cells.Item(13, 4).Value = 1.042314059082656
$ws.Cells.Item(13, 5).Value = 1.049644279964475
$ws.Cells.Item(13, 6).Value = 1.055471786393279
$ws.Cells.Item(13, 9).Value = 1.042002483313759
$ws.Cells.Item(13, 10).Value = 1.038474629385306
$ws.Cells.Item(13, 11).Value = 1.045893001005059
$ws.Cells.Item(13, 12).Value = 1.05319633651231
$ws.Cells.Item(13, 13).Value = 1.059002763345171
$ws.Cells.Item(13, 14).Value = 1.016800169140241

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.031945501364101
$ws.Cells.Item(14, 4).Value = 1.042416004433342
$ws.Cells.Item(14, 5).Value = 1.049783363678092
$ws.Cells.Item(14, 6).Value = 1.055602736866041
$ws.Cells.Item(14, 9).Value = 1.042036091572584
$ws.Cells.Item(14, 10).Value = 1.038542946779724
$ws.Cells.Item(14, 11).Value = 1.04596317920212
$ws.Cells.Item(14, 12).Value = 1.053303775910145
$ws.Cells.Item(14, 13).Value = 1.059102301120944
$ws.Cells.Item(14, 14).Value = 1.016822851025668

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032025167403501
$ws.Cells.Item(15, 4).Value = 1.04247882995864
$ws.Cells.Item(15, 5).Value = 1.049869092964988
$ws.Cells.Item(15, 6).Value = 1.055683446061424
$ws.Cells.Item(15, 9).Value = 1.042056774732513
$ws.Cells.Item(15, 10).Value = 1.038585037238509
$ws.Cells.Item(15, 11).Value = 1.046006414481692
$ws.Cells.Item(15, 12).Value = 1.053369991833933
$ws.Cells.Item(15, 13).Value = 1.059163639388454
$ws.Cells.Item(15, 14).Value = 1.016836825087005

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.032489092605905
$ws.Cells.Item(16, 4).Value = 1.042844708281179
$ws.Cells.Item(16, 5).Value = 1.050368610554492
$ws.Cells.Item(16, 6).Value = 1.056153608148876
$ws.Cells.Item(16, 9).Value = 1.042176792566457
$ws.Cells.Item(16, 10).Value = 1.038829988180486
$ws.Cells.Item(16, 11).Value = 1.046258003027901
$ws.Cells.Item(16, 12).Value = 1.053755685309412
$ws.Cells.Item(16, 13).Value = 1.059520802249792
$ws.Cells.Item(16, 14).Value = 1.01691814459049

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032780305488379
$ws.Cells.Item(17, 4).Value = 1.043074394433468
$ws.Cells.Item(17, 5).Value = 1.050682415072424
$ws.Cells.Item(17, 6).Value = 1.056448878605194
$ws.Cells.Item(17, 9).Value = 1.042251752529192
$ws.Cells.Item(17, 10).Value = 1.038983608652525
$ws.Cells.Item(17, 11).Value = 1.046415764331684
$ws.Cells.Item(17, 12).Value = 1.05399787294838
$ws.Cells.Item(17, 13).Value = 1.059744969624786
$ws.Cells.Item(17, 14).Value = 1.016969140070335

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.032950238229369
$ws.Cells.Item(18, 4).Value = 1.043208430831498
$ws.Cells.Item(18, 5).Value = 1.050865620716208
$ws.Cells.Item(18, 6).Value = 1.056621230352533
$ws.Cells.Item(18, 9).Value = 1.042295358334024
$ws.Cells.Item(18, 10).Value = 1.039073201353332
$ws.Cells.Item(18, 11).Value = 1.046507764072998
$ws.Cells.Item(18, 12).Value = 1.054139227235973
$ws.Cells.Item(18, 13).Value = 1.05987576840396
$ws.Cells.Item(18, 14).Value = 1.016998879640681

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.033008193257401
$ws.Cells.Item(19, 4).Value = 1.043254144633628
$ws.Cells.Item(19, 5).Value = 1.050928117665454
$ws.Cells.Item(19, 6).Value = 1.056680019060665
$ws.Cells.Item(19, 9).Value = 1.042310206901617
$ws.Cells.Item(19, 10).Value = 1.039103748154267
$ws.Cells.Item(19, 11).Value = 1.046539130212303
$ws.Cells.Item(19, 12).Value = 1.054187440661298
$ws.Cells.Item(19, 13).Value = 1.059920375110768
$ws.Cells.Item(19, 14).Value = 1.017009019164266

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.032749053517104
$ws.Cells.Item(20, 4).Value = 1.043049744612469
$ws.Cells.Item(20, 5).Value = 1.050648729343794
$ws.Cells.Item(20, 6).Value = 1.056417185883886
$ws.Cells.Item(20, 9).Value = 1.042243722137558
$ws.Cells.Item(20, 10).Value = 1.038967127811848
$ws.Cells.Item(20, 11).Value = 1.046398840070148
$ws.Cells.Item(20, 12).Value = 1.053971879149634
$ws.Cells.Item(20, 13).Value = 1.059720913842525
$ws.Cells.Item(20, 14).Value = 1.016963669275587

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031907429878897
$ws.Cells.Item(21, 4).Value = 1.042385981225418
$ws.Cells.Item(21, 5).Value = 1.049742399649678
$ws.Cells.Item(21, 6).Value = 1.05556416973539
$ws.Cells.Item(21, 9).Value = 1.04202619977566
$ws.Cells.Item(21, 10).Value = 1.03852282945381
$ws.Cells.Item(21, 11).Value = 1.045942514266241
$ws.Cells.Item(21, 12).Value = 1.053272133731193
$ws.Cells.Item(21, 13).Value = 1.059072987690654
$ws.Cells.Item(21, 14).Value = 1.016816171982573

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031378940692176
$ws.Cells.Item(22, 4).Value = 1.041969240341365
$ws.Cells.Item(22, 5).Value = 1.049174090190882
$ws.Cells.Item(22, 6).Value = 1.055028990569373
$ws.Cells.Item(22, 9).Value = 1.04188838631387
$ws.Cells.Item(22, 10).Value = 1.038243386808157
$ws.Cells.Item(22, 11).Value = 1.045655436001175
$ws.Cells.Item(22, 12).Value = 1.052833002565026
$ws.Cells.Item(22, 13).Value = 1.058666035415367
$ws.Cells.Item(22, 14).Value = 1.016723390675708

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03165903835722
$ws.Cells.Item(23, 4).Value = 1.042190105645233
$ws.Cells.Item(23, 5).Value = 1.04947521544309
$ws.Cells.Item(23, 6).Value = 1.055312589800487
$ws.Cells.Item(23, 9).Value = 1.041961542924634
$ws.Cells.Item(23, 10).Value = 1.038391533167873
$ws.Cells.Item(23, 11).Value = 1.045807637078229
$ws.Cells.Item(23, 12).Value = 1.053065715340548
$ws.Cells.Item(23, 13).Value = 1.0588817276608
$ws.Cells.Item(23, 14).Value = 1.016772579807829

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032763174711802
$ws.Cells.Item(24, 4).Value = 1.043060882606975
$ws.Cells.Item(24, 5).Value = 1.050663949953387
$ws.Cells.Item(24, 6).Value = 1.056431506073021
$ws.Cells.Item(24, 9).Value = 1.042247351087726
$ws.Cells.Item(24, 10).Value = 1.038974574830117
$ws.Cells.Item(24, 11).Value = 1.04640648747668
$ws.Cells.Item(24, 12).Value = 1.053983624349525
$ws.Cells.Item(24, 13).Value = 1.059731783473968
$ws.Cells.Item(24, 14).Value = 1.016966141308469

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.034047327591602
$ws.Cells.Item(25, 4).Value = 1.044073887928128
$ws.Cells.Item(25, 5).Value = 1.052050022866716
$ws.Cells.Item(25, 6).Value = 1.057734860563635
$ws.Cells.Item(25, 9).Value = 1.042574419752872
$ws.Cells.Item(25, 10).Value = 1.039650705439779
$ws.Cells.Item(25, 11).Value = 1.047100639762158
$ws.Cells.Item(25, 12).Value = 1.055052343650818
$ws.Cells.Item(25, 13).Value = 1.060720013711363
$ws.Cells.Item(25, 14).Value = 1.017190551789341
